# Update F-column ("想去人数" / interested-count) values per the diff, across 3 worksheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Worksheet 1 / sheet1.xml) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1376
$ws1.Range("F3").Value = 1924
$ws1.Range("F4").Value = 908
$ws1.Range("F7").Value = 688
$ws1.Range("F11").Value = 2515
$ws1.Range("F12").Value = 1611
$ws1.Range("F13").Value = 1559
$ws1.Range("F14").Value = 315
$ws1.Range("F15").Value = 253
$ws1.Range("F16").Value = 629
$ws1.Range("F17").Value = 808
$ws1.Range("F18").Value = 90
$ws1.Range("F19").Value = 320
$ws1.Range("F20").Value = 1093
$ws1.Range("F22").Value = 35
$ws1.Range("F24").Value = 5318
$ws1.Range("F25").Value = 224
$ws1.Range("F26").Value = 693
$ws1.Range("F27").Value = 92
$ws1.Range("F28").Value = 165
$ws1.Range("F29").Value = 144
$ws1.Range("F30").Value = 237
$ws1.Range("F32").Value = 37
$ws1.Range("F33").Value = 1052
$ws1.Range("F34").Value = 765
$ws1.Range("F38").Value = 405
$ws1.Range("F39").Value = 1118
$ws1.Range("F40").Value = 140
$ws1.Range("F42").Value = 181
$ws1.Range("F43").Value = 132
$ws1.Range("F44").Value = 74

# --- Sheet "演出" (Worksheet 2 / sheet2.xml) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 431
$ws2.Range("F6").Value = 11

# --- Sheet "全部类型" (Worksheet 4 / sheet4.xml) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1376
$ws4.Range("F4").Value = 1924
$ws4.Range("F5").Value = 908
$ws4.Range("F9").Value = 688
$ws4.Range("F11").Value = 11
$ws4.Range("F15").Value = 2515
$ws4.Range("F16").Value = 1611
$ws4.Range("F17").Value = 1559
$ws4.Range("F18").Value = 315
$ws4.Range("F19").Value = 253
$ws4.Range("F20").Value = 629
$ws4.Range("F22").Value = 808
$ws4.Range("F23").Value = 90
$ws4.Range("F24").Value = 320
$ws4.Range("F25").Value = 1093
$ws4.Range("F26").Value = 35
$ws4.Range("F28").Value = 5319
$ws4.Range("F29").Value = 224
$ws4.Range("F30").Value = 693
$ws4.Range("F31").Value = 92
$ws4.Range("F32").Value = 165
$ws4.Range("F33").Value = 144
$ws4.Range("F34").Value = 237
$ws4.Range("F36").Value = 37
$ws4.Range("F37").Value = 1052
$ws4.Range("F38").Value = 765
$ws4.Range("F40").Value = 405
$ws4.Range("F41").Value = 1118
$ws4.Range("F42").Value = 140
$ws4.Range("F44").Value = 181
$ws4.Range("F45").Value = 132
$ws4.Range("F46").Value = 74
